$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 (bold, bordered, centered) to the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
  @(3,4),
  @(8,9),
  @(8,8),
  @(7,8),
  @(9,9),
  @(5,6),
  @(4,6),
  @(9,9),
  @(8,9),
  @(4,6),
  @(5,5),
  @(5,5),
  @(8,8),
  @(8,8),
  @(8,8),
  @(5,6),
  @(7,7),
  @(8,8),
  @(8,8),
  @(6,6),
  @(9,9),
  @(9,9),
  @(6,6),
  @(6,6),
  @(7,7),
  @(7,7),
  @(7,7),
  @(6,6),
  @(7,7),
  @(4,4),
  @(8,8),
  @(7,7),
  @(5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).Value = $data[$i][0]
  $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
